# Update the degree abbreviations on the Author_form worksheet to include
# periods, e.g. "MSc Forestry" -> "M.Sc. Forestry" and
# "BSc Microbiology (Co-op)" -> "B.Sc. Microbiology (Co-op)"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A15").Value = "M.Sc. Forestry"
$ws.Range("A16").Value = "B.Sc. Microbiology (Co-op)"

# Update the saved selection to match the author's cursor position after edit
$ws.Range("A16").Select()
